$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.05
$ws.Range("G2").Value = 3.7
$ws.Range("H2").Value = 2.28
$ws.Range("I2").Value = 2.68
$ws.Range("J2").Value = 2.92
$ws.Range("K2").Value = 3.7
$ws.Range("N2").Value = 2.78
$ws.Range("Q2").Value = 1.98
$ws.Range("V2").Value = 1.6
$ws.Range("W2").Value = 1.37
$ws.Range("X2").Value = 14.5
$ws.Range("Y2").Value = 11
$ws.Range("Z2").Value = 18.5
$ws.Range("AA2").Value = 44
$ws.Range("AB2").Value = 14
$ws.Range("AE2").Value = 36
$ws.Range("AF2").Value = 28
$ws.Range("AN2").Value = 60
$ws.Range("AO2").Value = 32

# Row 3
$ws.Range("F3").Value = 1.47
$ws.Range("G3").Value = 1.69
$ws.Range("H3").Value = 6.2
$ws.Range("I3").Value = 10.5
$ws.Range("J3").Value = 3.65
$ws.Range("K3").Value = 5.5
$ws.Range("L3").Value = 1.37
$ws.Range("N3").Value = 3.15
$ws.Range("O3").Value = 1.3
$ws.Range("P3").Value = 1.86
$ws.Range("Q3").Value = 1.87
$ws.Range("R3").Value = 1.33
$ws.Range("S3").Value = 1.87
$ws.Range("T3").Value = 2
$ws.Range("U3").Value = 1.77
$ws.Range("V3").Value = 1.1
$ws.Range("W3").Value = 2.44
$ws.Range("AC3").Value = 12
$ws.Range("AF3").Value = 10.5
$ws.Range("AN3").Value = 11

# Row 4
$ws.Range("F4").Value = 1.27
$ws.Range("G4").Value = 1.37
$ws.Range("H4").Value = 6.6
$ws.Range("I4").Value = 15.5
$ws.Range("J4").Value = 4.9
$ws.Range("K4").Value = 7.6
$ws.Range("L4").Value = 1.27
$ws.Range("N4").Value = 3.7
$ws.Range("P4").Value = 2.04
$ws.Range("R4").Value = 1.44
$ws.Range("S4").Value = 2.5
$ws.Range("T4").Value = 2.14
$ws.Range("U4").Value = 1.66
$ws.Range("V4").Value = 1.07
$ws.Range("AD4").Value = 60
$ws.Range("AJ4").Value = 12
$ws.Range("AN4").Value = 6.4

# Row 5
$ws.Range("F5").Value = 2.68
$ws.Range("I5").Value = 2.78
$ws.Range("K5").Value = 4
$ws.Range("L5").Value = 1.33
$ws.Range("V5").Value = 1.56
$ws.Range("W5").Value = 1.5
$ws.Range("AE5").Value = 32
$ws.Range("AF5").Value = 22
$ws.Range("AI5").Value = 40
$ws.Range("AM5").Value = 90
$ws.Range("AN5").Value = 26

# Row 6
$ws.Range("N6").Value = 1.02

# Row 7
$ws.Range("F7").Value = 2.68
$ws.Range("G7").Value = 3.25
$ws.Range("H7").Value = 2.46
$ws.Range("I7").Value = 3.05
$ws.Range("L7").Value = 1.44
$ws.Range("N7").Value = 2.78
$ws.Range("O7").Value = 1.37
$ws.Range("P7").Value = 1.69
$ws.Range("R7").Value = 1.26
$ws.Range("T7").Value = 1.81
$ws.Range("U7").Value = 1.95
$ws.Range("V7").Value = 1.49
$ws.Range("W7").Value = 1.45
$ws.Range("AO7").Value = 42

# Row 8
$ws.Range("F8").Value = 1.92
$ws.Range("G8").Value = 2
$ws.Range("J8").Value = 3.6
$ws.Range("K8").Value = 3.85
$ws.Range("N8").Value = 3.5
$ws.Range("P8").Value = 1.86
$ws.Range("Q8").Value = 1.95
$ws.Range("S8").Value = 3.4
$ws.Range("T8").Value = 1.83
$ws.Range("U8").Value = 2
$ws.Range("W8").Value = 2
$ws.Range("X8").Value = 14.5
$ws.Range("Y8").Value = 16.5
$ws.Range("Z8").Value = 36
$ws.Range("AB8").Value = 9
$ws.Range("AC8").Value = 8.6
$ws.Range("AD8").Value = 19
$ws.Range("AF8").Value = 12
$ws.Range("AG8").Value = 11
$ws.Range("AH8").Value = 20
$ws.Range("AJ8").Value = 23
$ws.Range("AK8").Value = 22
$ws.Range("AL8").Value = 40
$ws.Range("AN8").Value = 14.5
$ws.Range("AO8").Value = 80

# Row 9
$ws.Range("K9").Value = 4
$ws.Range("L9").Value = 1.28
$ws.Range("X9").Value = 21
$ws.Range("Y9").Value = 15.5
$ws.Range("AE9").Value = 34
$ws.Range("AL9").Value = 40

# Row 10
$ws.Range("R10").Value = 1.71

# Row 11
$ws.Range("G11").Value = 2.9
$ws.Range("I11").Value = 2.82
$ws.Range("K11").Value = 3.7
$ws.Range("L11").Value = 1.42
$ws.Range("V11").Value = 1.55
$ws.Range("AH11").Value = 19.5
$ws.Range("AJ11").Value = 50
$ws.Range("AM11").Value = 1000
$ws.Range("AO11").Value = 28

# Row 12
$ws.Range("F12").Value = 1.81
$ws.Range("H12").Value = 3.6
$ws.Range("I12").Value = 4.9
$ws.Range("N12").Value = 3.6
$ws.Range("O12").Value = 1.27
$ws.Range("T12").Value = 1.75
$ws.Range("V12").Value = 1.26
$ws.Range("AN12").Value = 12.5

# Row 13
$ws.Range("F13").Value = 1.85
$ws.Range("G13").Value = 1.96
$ws.Range("U13").Value = 1.88
$ws.Range("X13").Value = 13.5
$ws.Range("Y13").Value = 15
$ws.Range("Z13").Value = 40
$ws.Range("AB13").Value = 8.800000000000001
$ws.Range("AC13").Value = 9.6
$ws.Range("AD13").Value = 22
$ws.Range("AE13").Value = 75
$ws.Range("AF13").Value = 12
$ws.Range("AG13").Value = 11
$ws.Range("AH13").Value = 26
$ws.Range("AN13").Value = 20
$ws.Range("AO13").Value = 110
